# Convert the opening "Heading1" title + bold "By <Author>" byline into a
# pandoc-style title block: a "Title"-styled paragraph (word-split into
# separate runs) followed by an "Authors"-styled paragraph (also word-split),
# with the old bookmarks around the heading removed.

$d = $word.ActiveDocument

# --- Step 1: drop the old heading paragraph (text + its trailing paragraph
# mark) that sits between the bookmarkStart/bookmarkEnd pair. Once the
# paragraph is gone, the bookmark markers collapse to the same (zero-width)
# position and can be deleted there.
$headingPara = $d.Paragraphs.Item(1)
$headingRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End)
$null = $headingRange.Delete()

# Removing a zero-length range sitting exactly on a bookmark marker deletes
# that marker; do it twice to clear both bookmarkStart and bookmarkEnd.
$null = $d.Range(0, 0).Delete()
$null = $d.Range(0, 0).Delete()

# --- Step 2: the byline paragraph ("By Dorothy Day") is now paragraph 1.
# Replace it (and implicitly insert the new title paragraph before it) with
# two freshly-built paragraphs, each word (and the spaces between words)
# split into its own run, mirroring the target markup exactly.
$byline = $d.Paragraphs.Item(1)
$target = $byline.Range

function New-RunsXml([string[]]$pieces) {
    $sb = ""
    foreach ($piece in $pieces) {
        $escaped = $piece -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
        $sb += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }
    return $sb
}

$titleWords = @("What", " ", "is", " ", "Happening", "?", " ", "Trial", " ", "Continued", " ", "Until", " ", "Nov", ".", " ", "16")
$authorWords = @("Dorothy", " ", "Day")

$titleRuns = New-RunsXml $titleWords
$authorRuns = New-RunsXml $authorWords

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' `
  + '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>' `
  + '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>' `
  + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $target.InsertXML($xml)
